# ----------------------------------------------------------------------------
# Adds two new spec columns (CU, CV) to Sheet1 and two new simulation-case
# rows (7, 8) that carry values for every existing column plus the two new
# ones. Mirrors the source edit: new headers in CU1:CV1, blank placeholder
# cells for the two new columns on the pre-existing data rows (2-6), and two
# brand-new fully populated rows at the bottom.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells CU1 / CV1 -------------------------------------------
$ws.Range("CU1").Value = "Tatoray Stripper C620 Operation_Specifications_Spec 2 : Distillate Rate_m3/hr"
$ws.Range("CV1").Value = "Benzene Column C660 Operation_Specifications_Spec 3 : Toluene in Benzene_ppmw"

# Match the bold / centered / bordered header style used by the rest of row 1
# (copy the formatting from the last existing header cell, CT1).
$ws.Range("CT1").Copy() | Out-Null
$ws.Range("CU1:CV1").PasteSpecial(-4122) | Out-Null

# --- Extend the existing data rows (2-6) with blank CU/CV cells -----------
# These rows don't have values for the two new spec columns yet, but the
# column range still needs to cover them, so materialize empty cells.
$ws.Range("CU2:CV6").Borders.LineStyle = -4142

# --- Shared payload for the two new simulation rows ------------------------
# Columns B..CV, in order (99 values); column A is set separately per row.
$rowValues = @(
    163,
    1.116022501559,
    11.15222454,
    48.74597931,
    15,
    2.3389766495,
    74.30225372,
    17.46782494,
    72,
    0.09109998293844999,
    20.85439682,
    74.68977356000001,
    70,
    980,
    36,
    0,
    0.01829845831,
    0.4437376857,
    0.4127115904999999,
    0.1285301149,
    0.07167375832999999,
    0.01489215251,
    16.52599144,
    0.007658495568,
    0.01510457043,
    0,
    0.01225359179,
    46.08435059,
    0.008169061504,
    0.009615667164000001,
    0.009615667164000001,
    0.7336538434000001,
    5.038924217,
    10.67489147,
    4.729251385,
    0.01395548135,
    0.01189399883,
    0.002637926256,
    0.01463876944,
    0.9881169795999999,
    0.5178464054999999,
    2.787770271,
    0.2909455299,
    7.119931221,
    0,
    1.136334419,
    0.05123569444,
    0.0009149230900000001,
    0.020128306,
    0.7639607191,
    0,
    0.8801559806,
    0.3220529258,
    0.1381533742,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0.00009999997565000002,
    0.0004999999655,
    20.85439682,
    0.000299999956,
    0.0005999999703,
    0,
    0.002599999309,
    74.68977356000001,
    0.01079999842,
    0.009599998593,
    0.007899997756000001,
    1.97609961,
    0.9054998159000001,
    1.393299699,
    0.08289998025,
    0.01729999669,
    0,
    0.03829999268,
    0,
    0,
    0,
    0,
    0,
    0,
    0.003099999623,
    0,
    0,
    0.005499999039,
    0.001399999601,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0.009999999775999999,
    9.999955072999999
)

# --- New row 7 ("a1") --------------------------------------------------
$ws.Range("A7").Value = "a1"
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(7, $i + 2).Value = $rowValues[$i]
}

# --- New row 8 ("b2") --------------------------------------------------
$ws.Range("A8").Value = "b2"
for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(8, $i + 2).Value = $rowValues[$i]
}

# Column A on the new rows uses the same style as the other label cells
# (e.g. A5), so copy that formatting over too.
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A7:A8").PasteSpecial(-4122) | Out-Null

Write-Output "edit complete"
